$d = $word.ActiveDocument

function Insert-XmlRuns($rng, $innerXml) {
    $pkgXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkgXml)
}

# --- Edit 1: "housing price" source-file paragraph ---
$rng1 = $d.Content
$rng1.Find.Execute("(it is downloaded and in data folder [name starts with zip_zhvi_uc_sfrc]", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$inner1 = '<w:r><w:t>(</w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>it</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> is downloaded and in data folder [name starts with </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Metro_mlp_uc_sfrcondo_sm_month</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>]</w:t></w:r>'

Insert-XmlRuns $rng1 $inner1
$rng1.Delete()

# --- Edit 2: "Using RegionName (zip codes) ... " merging-steps bullet ---
$rng2 = $d.Content
$rng2.Find.Execute(" RegionName (zip codes) to merge this data with the demographics data . ", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$inner2 = '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>RegionName</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> (zip codes)</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> or </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>StateName</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> to merge this data with the demographics data </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">. </w:t></w:r>'

Insert-XmlRuns $rng2 $inner2
$rng2.Delete()
